# Disable smart-quote autocorrect so straight apostrophes are preserved
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoCorrectAsYouType = $false } catch {}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title.
#    We copy the existing bottom "bold title" paragraph (same run layout:
#    empty run + bold run) so the new paragraph gets the identical
#    <w:r/><w:r><w:rPr><w:b/></w:rPr>... structure, then fix the text.
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($paraCount - 1)
$boldTitlePara.Range.Copy()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Replace the pasted (bold) text with "Meta description"
$metaParaStart = $metaPara.Range.Start
$metaParaEnd = $metaPara.Range.End
$boldTextRange = $d.Range($metaParaStart, $metaParaEnd - 1)
$boldTextRange.Text = "Meta description"

# Append the rest of the sentence as a non-bold run right after it
$afterBoldStart = $metaParaStart + ("Meta description").Length
$afterBoldRange = $d.Range($afterBoldStart, $afterBoldStart)
$afterBoldRange.InsertAfter(": Join Cornelius the cat and play this NetEnt slot game for free. Featuring Cash Drop symbol, Free Spins bonus, and unique symbols. Available on desktop and mobile.")
$afterBoldRange.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Remove the duplicate bold title paragraph at the bottom of the
#    document (it has now been relocated to the top as "Meta description").
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$bottomBoldPara = $d.Paragraphs.Item($paraCount - 1)
$bottomBoldPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-prompt copy, keeping its italic formatting intact.
# ---------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($paraCount)
$italicTextRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End - 1)
$italicTextRange.Text = "Create a cartoon-style feature image showcasing Cornelius, the gluttonous kitty with glasses. The image should feature a happy Maya warrior with glasses, who is enjoying a sweet treat with Cornelius. The Maya warrior should be dressed in traditional warrior attire with a headdress, while Cornelius should be sitting on a giant cookie, wearing his green shirt that's too tight for him. The background should be a mint green kitchen with a fridge on the right and Cornelius' red bowl on the left. The image should have bright and cheerful colors to reflect the fun theme of the game."
